$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two name/link swaps) per latest scrape.
# Values are written with a leading apostrophe to force text storage (matching the
# original inlineStr cell type), then the style is reset to Normal so no stray
# number-format/style is left applied to the cell.

# Row 2
$ws.Range("D2").Value = "'62.739.88"
$ws.Range("E2").Value = "'  -5.82%  "
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.063.02"
$ws.Range("E3").Value = "'  -6.27%  "
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'541.54"
$ws.Range("E5").Value = "'  -7.39%  "
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'133.86"
$ws.Range("E6").Value = "'  -12.86%  "
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'3.063.32"
$ws.Range("E8").Value = "'  -6.20%  "
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "'  -5.39%  "
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.155"
$ws.Range("E10").Value = "'  -6.65%  "
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'6.12"
$ws.Range("E11").Value = "'  -13.65%  "
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "'  -5.92%  "
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'34.84"
$ws.Range("E13").Value = "'  -8.70%  "
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.0000221"
$ws.Range("E14").Value = "'  -6.83%  "
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'3.547.30"
$ws.Range("E15").Value = "'  -6.44%  "
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'62.621.09"
$ws.Range("E16").Value = "'  -6.10%  "
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "'  -3.58%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.050.50"
$ws.Range("E18").Value = "'  -6.56%  "
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'6.69"
$ws.Range("E19").Value = "'  -6.99%  "
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'481.04"
$ws.Range("E20").Value = "'  -13.89%  "
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'13.46"
$ws.Range("E21").Value = "'  -8.02%  "
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.711"
$ws.Range("E22").Value = "'  -5.24%  "
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'7.21"
$ws.Range("E23").Value = "'  -8.84%  "
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'78.72"
$ws.Range("E24").Value = "'  -4.18%  "
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'12.09"
$ws.Range("E25").Value = "'  -11.54%  "
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  +0.37%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Value = "'RenderToken"
$ws.Range("C27").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'8.33"
$ws.Range("E27").Value = "'  -10.55%  "
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = "'PancakeSwap"
$ws.Range("C28").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.72"
$ws.Range("E28").Value = "'  -9.57%  "
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "'  -0.21%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'1.93"
$ws.Range("E30").Value = "'  -15.23%  "
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'26.30"
$ws.Range("E31").Value = "'  -6.12%  "
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  -7.03%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("B33").Value = "'Stacks"
$ws.Range("C33").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.44"
$ws.Range("E33").Value = "'  -12.64%  "
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").Value = "'OKB"
$ws.Range("C34").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'58.69"
$ws.Range("E34").Value = "'  +5.75%  "
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'493.37"
$ws.Range("E35").Value = "'  -13.65%  "
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'5.97"
$ws.Range("E36").Value = "'  -7.42%  "
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'5.11"
$ws.Range("E37").Value = "'  -11.04%  "
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'3.136.29"
$ws.Range("E38").Value = "'  -2.79%  "
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.0395"
$ws.Range("E39").Value = "'  -13.53%  "
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.0798"
$ws.Range("E40").Value = "'  -8.54%  "
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "'  -11.61%  "
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'8.08"
$ws.Range("E42").Value = "'  -7.34%  "
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'2.56"
$ws.Range("E43").Value = "'  -15.95%  "
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.253"
$ws.Range("E44").Value = "'  -10.97%  "
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  +0.04%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'2.04"
$ws.Range("E46").Value = "'  -12.25%  "
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'24.83"
$ws.Range("E47").Value = "'  -7.47%  "
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.108"
$ws.Range("E48").Value = "'  -5.24%  "
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'117.81"
$ws.Range("E49").Value = "'  -6.40%  "
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.0₃0513"
$ws.Range("E50").Value = "'  -8.84%  "
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("B51").Value = "'ThetaToken"
$ws.Range("C51").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.02"
$ws.Range("E51").Value = "'  -9.67%  "
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
